$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record is inserted at the top of the "Arveja Verde" history
# (row 113). Every existing record from row 113 through row 141 is pushed down
# one row (so the old row 141 becomes the new row 142), and the brand-new
# record's data is written into row 113.

# Columns that vary between records.
$varCols = @("D", "H", "J", "K", "L", "M", "N", "O", "P")

# Columns that are constant for every record in this sheet/subset.
$fixedCols = @("A", "B", "C", "E", "F", "G", "I", "Q", "R")

# Shift existing records down by one row, working from the bottom up so that
# source data is never overwritten before it's copied (row 141 -> 142,
# row 140 -> 141, ... , row 113 -> 114).
for ($r = 141; $r -ge 113; $r--) {
    foreach ($c in $varCols) {
        $ws.Range("$c$($r + 1)").Value2 = $ws.Range("$c$r").Value2
    }
}

# Row 142 is a brand-new row, so also copy over the constant columns and the
# date cell's number format (so it keeps the same date style as the rest of
# column D).
foreach ($c in $fixedCols) {
    $ws.Range("$c" + "142").Value2 = $ws.Range("$c" + "141").Value2
}
$ws.Range("D142").NumberFormat = $ws.Range("D141").NumberFormat

# Finally, write the new record's data into row 113 (H113/N113/O113 keep the
# values already shifted into place: "Perfection" / "$/malla 25 kilos" /
# "Provincia de Huasco").
$ws.Range("D113").Value2 = 44841
$ws.Range("J113").Value2 = 70
$ws.Range("K113").Value2 = 34000
$ws.Range("L113").Value2 = 34000
$ws.Range("M113").Value2 = 34000
$ws.Range("P113").Value2 = 1360
